$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'288.45"
$ws.Range("E2").Value = "'-4.10%"

$ws.Range("D3").Value = "'30.38"
$ws.Range("E3").Value = "'-6.38%"

$ws.Range("D4").Value = "'4.952"
$ws.Range("E4").Value = "'-2.27%"

$ws.Range("D5").Value = "'0.07195"
$ws.Range("E5").Value = "'-6.00%"

$ws.Range("D6").Value = "'1.786"
$ws.Range("E6").Value = "'-16.00%"

$ws.Range("D7").Value = "'7.562"
$ws.Range("E7").Value = "'-3.69%"

$ws.Range("D8").Value = "'3.719"
$ws.Range("E8").Value = "'-1.75%"

$ws.Range("D9").Value = "'0.8994"
$ws.Range("E9").Value = "'-2.32%"

$ws.Range("D10").Value = "'0.1663"
$ws.Range("E10").Value = "'-5.51%"

$ws.Range("E11").Value = "'-1.56%"

$ws.Range("D12").Value = "'0.07955"
$ws.Range("E12").Value = "'-6.13%"

$ws.Range("D13").Value = "'0.03031"
$ws.Range("E13").Value = "'-0.95%"

$ws.Range("D14").Value = "'0.1001"
$ws.Range("E14").Value = "'0.05%"

$ws.Range("D15").Value = "'0.001495"
$ws.Range("E15").Value = "'-2.16%"

$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005750"
$ws.Range("E16").Value = "'0.18%"

$ws.Range("B17").Value = "UpBots"
$ws.Range("C17").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D17").Value = "'0.007492"
$ws.Range("E17").Value = "'-0.08%"

$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.474"
$ws.Range("E18").Value = "'0.35%"

$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.068"
$ws.Range("E19").Value = "'-3.88%"

$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3302"
$ws.Range("E20").Value = "'-1.19%"

$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1294"
$ws.Range("E21").Value = "'-1.73%"

$ws.Range("B22").Value = "MCDex"
$ws.Range("C22").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D22").Value = "'3.964"
$ws.Range("E22").Value = "'-7.16%"

$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").Value = "'0.2100"
$ws.Range("E23").Value = "'6.15%"

$ws.Range("B24").Value = "CoinExToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D24").Value = "'0.04505"
$ws.Range("E24").Value = "'-0.46%"

$ws.Range("E25").Value = "'-2.02%"

$ws.Range("D26").Value = "'0.004631"
$ws.Range("E26").Value = "'-3.94%"

$ws.Range("D27").Value = "'0.0001300"
$ws.Range("E27").Value = "'3.82%"

$ws.Range("D39").Value = "'0.01563"
$ws.Range("E39").Value = "'-8.14%"

$ws.Range("D40").Value = "'0.04322"
$ws.Range("E40").Value = "'-7.53%"

$ws.Range("D41").Value = "'0.007332"
$ws.Range("E41").Value = "'-1.68%"

$ws.Range("D43").Value = "'0.1304"
$ws.Range("E43").Value = "'-3.55%"

$ws.Range("E44").Value = "'-13.82%"

$ws.Range("D45").Value = "'0.009423"
$ws.Range("E45").Value = "'-11.36%"

$ws.Range("D46").Value = "'0.00005898"
$ws.Range("E46").Value = "'-5.21%"

$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.16%"

$ws.Range("D48").Value = "'2.259"
$ws.Range("E48").Value = "'115.87%"

$ws.Range("E49").Value = "'-0.16%"

$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.16%"

$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.16%"
